$wb = $excel.ActiveWorkbook

# --- Hoja1: update Ref labels (col A) and Cantidad values (col B) for rows 2-12 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$ws1.Range("A2").Value = "RH150010"
$ws1.Range("B2").Value = 163.37249999999997

$ws1.Range("A3").Value = "G4506"
$ws1.Range("B3").Value = 1

$ws1.Range("A4").Value = "E2671"
$ws1.Range("B4").Value = 1

$ws1.Range("A5").Value = "AIN0L11"
$ws1.Range("B5").Value = 1

$ws1.Range("A6").Value = "M0F015"
$ws1.Range("B6").Value = 44

$ws1.Range("A7").Value = "T1502015"
$ws1.Range("B7").Value = 56

$ws1.Range("A8").Value = "TS03981"
$ws1.Range("B8").Value = 92

$ws1.Range("A9").Value = "BSCO0051"
$ws1.Range("B9").Value = 46

$ws1.Range("A10").Value = "SP1"
$ws1.Range("B10").Value = 112

$ws1.Range("A11").Value = "MG023"
$ws1.Range("B11").Value = 14.904999999999998

$ws1.Range("A12").Value = "#02"
$ws1.Range("B12").Value = 1

# --- Hoja2: update quantity value ---
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Range("A2").Value = 3
